$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update (values scraped on Thu Sep 14 10:49:53 UTC 2023).
# D-column cells whose new text parses as a plain number must be forced to
# stay text (matching the source inlineStr cells) without leaving a stray
# number-format style behind, so: mark as Text, assign, then ClearFormats()
# to drop back to the default (unstyled) cell xf.

$ws.Range('D2').Value = '26.349.60'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.621.86'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  +0.11%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.0617'
$c.ClearFormats()
$ws.Range('E9').Value = '  +0.32%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.ClearFormats()
$ws.Range('E10').Value = '  +3.52%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0813'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.847.03'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '1.640.72'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '26.354.69'
$ws.Range('E16').Value = '  +0.65%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '62.52'
$c.ClearFormats()
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  -0.14%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '202.21'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.00%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.07%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '9.30'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.30%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '6.06'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('E24').Value = '  -3.24%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '144.45'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -2.03%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '15.19'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.13%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '6.59'
$c.ClearFormats()
$ws.Range('E29').Value = '  +1.07%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.0516'
$c.ClearFormats()
$ws.Range('E30').Value = '  +5.90%  '
$ws.Range('E31').Value = '  +0.39%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.20'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('E35').Value = '  +2.15%  '
$ws.Range('D36').Value = '1.162.85'
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('E37').Value = '  +0.04%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.802'
$c.ClearFormats()
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  -0.08%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '5.39'
$c.ClearFormats()
$ws.Range('E42').Value = '  +4.02%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.782'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '1.759.37'
$ws.Range('E44').Value = '  +1.24%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '92.16'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  +9.70%  '
$ws.Range('E47').Value = '  +0.96%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '53.91'
$c.ClearFormats()
$ws.Range('E48').Value = '  -0.31%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0508'
$c.ClearFormats()
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('E51').Value = '  -0.31%  '
